# Fix bug with hours:
# - G2/H2 held the wrong labels ("HDJ 1"/"HDJ 2") instead of the actual
#   afternoon time slots ("14h15"/"15h00").
# - Column A (rows 3-17) had stray room/HDJ labels left over from an
#   earlier layout; clear them out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "14h15"
$ws.Range("H2").Value = "15h00"

$ws.Range("A3:A17").ClearContents()

$ws.Range("A3:A17").Select()
